$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.498.17"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.866.42"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07820"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9907"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").Value = "1.874.73"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.910"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.633"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06931"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009936"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "28.542.33"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.250"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.081"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.102.43"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.679"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.880"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9026"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.266"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.263"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02036"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.602"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.635"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07135"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.134"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.801"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
